$wb = $excel.ActiveWorkbook

$changes = @(
    @{Sheet='Citywide Totals'; Cell='K2'; Value=3083},
    @{Sheet='Citywide Totals'; Cell='K3'; Value=3040},
    @{Sheet='Citywide Totals'; Cell='B4'; Value=1700},
    @{Sheet='Citywide Totals'; Cell='H4'; Value=1728},
    @{Sheet='Citywide Totals'; Cell='J4'; Value=1818},
    @{Sheet='Citywide Totals'; Cell='K4'; Value=619},
    @{Sheet='Citywide Totals'; Cell='K5'; Value=200},
    @{Sheet='Citywide Totals'; Cell='K6'; Value=3611},
    @{Sheet='Citywide Totals'; Cell='B7'; Value=23333},
    @{Sheet='Citywide Totals'; Cell='H7'; Value=26041},
    @{Sheet='Citywide Totals'; Cell='J7'; Value=29288},
    @{Sheet='Citywide Totals'; Cell='K7'; Value=10553},
    @{Sheet='Logan Square'; Cell='K3'; Value=31},
    @{Sheet='Logan Square'; Cell='K7'; Value=145},
    @{Sheet='Austin'; Cell='K2'; Value=206},
    @{Sheet='Austin'; Cell='K3'; Value=210},
    @{Sheet='Austin'; Cell='K4'; Value=38},
    @{Sheet='Austin'; Cell='K5'; Value=18},
    @{Sheet='Austin'; Cell='K6'; Value=229},
    @{Sheet='Austin'; Cell='K7'; Value=701},
    @{Sheet='South Chicago'; Cell='K2'; Value=87},
    @{Sheet='South Chicago'; Cell='K4'; Value=13},
    @{Sheet='Garfield Park'; Cell='K2'; Value=117},
    @{Sheet='Garfield Park'; Cell='K3'; Value=151},
    @{Sheet='Garfield Park'; Cell='K7'; Value=417},
    @{Sheet='West Pullman'; Cell='K3'; Value=63},
    @{Sheet='West Pullman'; Cell='K4'; Value=10},
    @{Sheet='West Pullman'; Cell='K7'; Value=174},
    @{Sheet='Grand Crossing'; Cell='K2'; Value=94},
    @{Sheet='Grand Crossing'; Cell='K3'; Value=121},
    @{Sheet='Grand Crossing'; Cell='K6'; Value=111},
    @{Sheet='Grand Crossing'; Cell='K7'; Value=355},
    @{Sheet='Woodlawn'; Cell='K6'; Value=50},
    @{Sheet='Woodlawn'; Cell='K7'; Value=190},
    @{Sheet='By Neighborhood'; Cell='K4'; Value=35},
    @{Sheet='By Neighborhood'; Cell='K7'; Value=310},
    @{Sheet='By Neighborhood'; Cell='K8'; Value=701},
    @{Sheet='By Neighborhood'; Cell='K11'; Value=225},
    @{Sheet='By Neighborhood'; Cell='K15'; Value=109},
    @{Sheet='By Neighborhood'; Cell='K18'; Value=75},
    @{Sheet='By Neighborhood'; Cell='K19'; Value=320},
    @{Sheet='By Neighborhood'; Cell='K20'; Value=246},
    @{Sheet='By Neighborhood'; Cell='K25'; Value=42},
    @{Sheet='By Neighborhood'; Cell='K27'; Value=108},
    @{Sheet='By Neighborhood'; Cell='K29'; Value=554},
    @{Sheet='By Neighborhood'; Cell='K33'; Value=417},
    @{Sheet='By Neighborhood'; Cell='K36'; Value=124},
    @{Sheet='By Neighborhood'; Cell='K37'; Value=355},
    @{Sheet='By Neighborhood'; Cell='K42'; Value=364},
    @{Sheet='By Neighborhood'; Cell='K43'; Value=94},
    @{Sheet='By Neighborhood'; Cell='K48'; Value=127},
    @{Sheet='By Neighborhood'; Cell='K51'; Value=116},
    @{Sheet='By Neighborhood'; Cell='K52'; Value=290},
    @{Sheet='By Neighborhood'; Cell='K53'; Value=145},
    @{Sheet='By Neighborhood'; Cell='K54'; Value=205},
    @{Sheet='By Neighborhood'; Cell='K60'; Value=63},
    @{Sheet='By Neighborhood'; Cell='B63'; Value=405},
    @{Sheet='By Neighborhood'; Cell='H63'; Value=281},
    @{Sheet='By Neighborhood'; Cell='J63'; Value=102},
    @{Sheet='By Neighborhood'; Cell='K64'; Value=65},
    @{Sheet='By Neighborhood'; Cell='K67'; Value=414},
    @{Sheet='By Neighborhood'; Cell='K71'; Value=32},
    @{Sheet='By Neighborhood'; Cell='K72'; Value=51},
    @{Sheet='By Neighborhood'; Cell='K73'; Value=96},
    @{Sheet='By Neighborhood'; Cell='K77'; Value=74},
    @{Sheet='By Neighborhood'; Cell='K78'; Value=137},
    @{Sheet='By Neighborhood'; Cell='K79'; Value=274},
    @{Sheet='By Neighborhood'; Cell='K82'; Value=12},
    @{Sheet='By Neighborhood'; Cell='K84'; Value=75},
    @{Sheet='By Neighborhood'; Cell='K85'; Value=503},
    @{Sheet='By Neighborhood'; Cell='K88'; Value=118},
    @{Sheet='By Neighborhood'; Cell='K89'; Value=140},
    @{Sheet='By Neighborhood'; Cell='K90'; Value=93},
    @{Sheet='By Neighborhood'; Cell='K91'; Value=110},
    @{Sheet='By Neighborhood'; Cell='K93'; Value=41},
    @{Sheet='By Neighborhood'; Cell='K94'; Value=130},
    @{Sheet='By Neighborhood'; Cell='K95'; Value=174},
    @{Sheet='By Neighborhood'; Cell='K97'; Value=91},
    @{Sheet='By Neighborhood'; Cell='K99'; Value=190},
    @{Sheet='By Neighborhood'; Cell='B101'; Value=23333},
    @{Sheet='By Neighborhood'; Cell='H101'; Value=26041},
    @{Sheet='By Neighborhood'; Cell='J101'; Value=29288},
    @{Sheet='By Neighborhood'; Cell='K101'; Value=10553},
    @{Sheet='North Lawndale'; Cell='K3'; Value=135},
    @{Sheet='North Lawndale'; Cell='K7'; Value=414},
    @{Sheet='South Deering'; Cell='K2'; Value=24},
    @{Sheet='South Deering'; Cell='K3'; Value=27},
    @{Sheet='South Deering'; Cell='K7'; Value=75},
    @{Sheet='Loop'; Cell='K3'; Value=62},
    @{Sheet='Loop'; Cell='K6'; Value=93},
    @{Sheet='Loop'; Cell='K7'; Value=205},
    @{Sheet='Englewood'; Cell='K2'; Value=150},
    @{Sheet='Englewood'; Cell='K3'; Value=193},
    @{Sheet='Englewood'; Cell='K6'; Value=169},
    @{Sheet='Englewood'; Cell='K7'; Value=554},
    @{Sheet='Lake View'; Cell='K2'; Value=20},
    @{Sheet='Lake View'; Cell='K3'; Value=26},
    @{Sheet='Lake View'; Cell='K4'; Value=17},
    @{Sheet='Lake View'; Cell='K6'; Value=64},
    @{Sheet='Lake View'; Cell='K7'; Value=127},
    @{Sheet='Chatham'; Cell='K2'; Value=107},
    @{Sheet='Chatham'; Cell='K3'; Value=83},
    @{Sheet='Chatham'; Cell='K6'; Value=103},
    @{Sheet='Chatham'; Cell='K7'; Value=320},
    @{Sheet='Humboldt Park'; Cell='K3'; Value=116},
    @{Sheet='Humboldt Park'; Cell='K4'; Value=15},
    @{Sheet='Humboldt Park'; Cell='K6'; Value=139},
    @{Sheet='Humboldt Park'; Cell='K7'; Value=364},
    @{Sheet='Rogers Park'; Cell='K2'; Value=40},
    @{Sheet='Rogers Park'; Cell='K7'; Value=137},
    @{Sheet='Washington Park'; Cell='K3'; Value=51},
    @{Sheet='Washington Park'; Cell='K7'; Value=110},
    @{Sheet='Roseland'; Cell='K2'; Value=95},
    @{Sheet='Roseland'; Cell='K3'; Value=97},
    @{Sheet='Roseland'; Cell='K7'; Value=274},
    @{Sheet='Near South Side'; Cell='K6'; Value=22},
    @{Sheet='Near South Side'; Cell='K7'; Value=65},
    @{Sheet='Chicago Lawn'; Cell='K4'; Value=8},
    @{Sheet='Chicago Lawn'; Cell='K7'; Value=246},
    @{Sheet='Calumet Heights'; Cell='K3'; Value=24},
    @{Sheet='Calumet Heights'; Cell='K7'; Value=75},
    @{Sheet='Grand Boulevard'; Cell='K2'; Value=48},
    @{Sheet='Grand Boulevard'; Cell='K4'; Value=12},
    @{Sheet='Grand Boulevard'; Cell='K7'; Value=124},
    @{Sheet='West Lawn'; Cell='K2'; Value=15},
    @{Sheet='West Lawn'; Cell='K7'; Value=41},
    @{Sheet='Auburn Gresham'; Cell='K2'; Value=106},
    @{Sheet='Auburn Gresham'; Cell='K3'; Value=96},
    @{Sheet='Auburn Gresham'; Cell='K6'; Value=83},
    @{Sheet='Auburn Gresham'; Cell='K7'; Value=310},
    @{Sheet='West Loop'; Cell='K2'; Value=36},
    @{Sheet='West Loop'; Cell='K3'; Value=25},
    @{Sheet='West Loop'; Cell='K6'; Value=56},
    @{Sheet='West Loop'; Cell='K7'; Value=130},
    @{Sheet='East Side'; Cell='K2'; Value=16},
    @{Sheet='East Side'; Cell='K7'; Value=42},
    @{Sheet='Brighton Park'; Cell='K2'; Value=37},
    @{Sheet='Brighton Park'; Cell='K6'; Value=36},
    @{Sheet='Brighton Park'; Cell='K7'; Value=109},
    @{Sheet='Belmont Cragin'; Cell='K2'; Value=67},
    @{Sheet='Belmont Cragin'; Cell='K6'; Value=85},
    @{Sheet='Belmont Cragin'; Cell='K7'; Value=225},
    @{Sheet='Portage Park'; Cell='K2'; Value=29},
    @{Sheet='Portage Park'; Cell='K7'; Value=96},
    @{Sheet='West Town'; Cell='K2'; Value=19},
    @{Sheet='West Town'; Cell='K4'; Value=5},
    @{Sheet='West Town'; Cell='K6'; Value=56},
    @{Sheet='West Town'; Cell='K7'; Value=91},
    @{Sheet='United Center'; Cell='K6'; Value=58},
    @{Sheet='United Center'; Cell='K7'; Value=118},
    @{Sheet='Uptown'; Cell='K2'; Value=33},
    @{Sheet='Uptown'; Cell='K6'; Value=40},
    @{Sheet='Uptown'; Cell='K7'; Value=140},
    @{Sheet='Edgewater'; Cell='K6'; Value=43},
    @{Sheet='Edgewater'; Cell='K7'; Value=108},
    @{Sheet='Washington Heights'; Cell='K2'; Value=35},
    @{Sheet='Washington Heights'; Cell='K7'; Value=93},
    @{Sheet='Little Italy, UIC'; Cell='K2'; Value=31},
    @{Sheet='Little Italy, UIC'; Cell='K6'; Value=42},
    @{Sheet='Little Italy, UIC'; Cell='K7'; Value=116},
    @{Sheet='Morgan Park'; Cell='K6'; Value=18},
    @{Sheet='Morgan Park'; Cell='K7'; Value=63},
    @{Sheet='Hyde Park'; Cell='K3'; Value=28},
    @{Sheet='Hyde Park'; Cell='K4'; Value=9},
    @{Sheet='Hyde Park'; Cell='K7'; Value=94},
    @{Sheet='South Shore'; Cell='K3'; Value=172},
    @{Sheet='South Shore'; Cell='K6'; Value=113},
    @{Sheet='South Shore'; Cell='K7'; Value=503},
    @{Sheet='Oakland'; Cell='K6'; Value=8},
    @{Sheet='Oakland'; Cell='K7'; Value=32},
    @{Sheet='Old Town'; Cell='K3'; Value=15},
    @{Sheet='Old Town'; Cell='K7'; Value=51},
    @{Sheet='Sheffield & DePaul'; Cell='K3'; Value=2},
    @{Sheet='Sheffield & DePaul'; Cell='K6'; Value=12},
    @{Sheet='Riverdale'; Cell='K3'; Value=27},
    @{Sheet='Riverdale'; Cell='K7'; Value=74},
    @{Sheet='Little Village'; Cell='K2'; Value=77},
    @{Sheet='Little Village'; Cell='K3'; Value=75},
    @{Sheet='Little Village'; Cell='K6'; Value=117},
    @{Sheet='Little Village'; Cell='K7'; Value=290},
    @{Sheet='Archer Heights'; Cell='K6'; Value=13},
    @{Sheet='Archer Heights'; Cell='K7'; Value=35}
)

foreach ($change in $changes) {
    $ws = $wb.Worksheets.Item($change.Sheet)
    $ws.Range($change.Cell).Value = $change.Value
}
